# "Removed address book from undo/redo"
#
# The original deck illustrates undo/redo state transitions using three
# "AddressBook" state boxes (ab0/ab1/ab2). This edit renames them to the
# "HotelManagementSystem" (hms0/hms1/hms2) example and repositions/resizes
# several shapes to make room for the longer label text.
#
# EMU <-> point conversion: PowerPoint's COM surface (Shape.Left/Top/Width/
# Height, Table Row.Height / Column.Width, ...) works in points, while the
# underlying OOXML stores English Metric Units (1 pt = 12700 EMU). A tiny
# epsilon is added before converting to points so that float round-tripping
# lands on the exact target EMU value instead of off-by-one.
function EMU([double]$emu) {
    return ($emu / 12700.0) + 0.00003
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1) "Down Arrow 49" - shift right/up slightly.
# ---------------------------------------------------------------------
$downArrow = $s.Shapes.Item(1)
$downArrow.Left = EMU 8070256
$downArrow.Top  = EMU 2135019

# ---------------------------------------------------------------------
# 2) "TextBox 15" (currentStatePointer = 1, top row) - shift right.
# ---------------------------------------------------------------------
$textBox15 = $s.Shapes.Item(2)
$textBox15.Left = EMU 3909780
$textBox15.Top  = EMU 2731691

# ---------------------------------------------------------------------
# Helper: rebuild a single-cell, single-row table's cell text.
# The COM text-range implementation for table cells can only overwrite
# the first run of a paragraph (it does not support deleting/replacing a
# second run), so the row is dropped and re-added to get a clean,
# single-run paragraph with the fully-updated label text.
# ---------------------------------------------------------------------
function Set-StateTableLabel($shape, [string]$text, [double]$widthEmu, [double]$heightEmu) {
    $tbl = $shape.Table
    $tbl.Rows.Item(1).Delete()
    [void]$tbl.Rows.Add()
    $cell = $tbl.Cell(1, 1)
    $tr = $cell.Shape.TextFrame.TextRange
    $tr.Text = $text
    $tr.Font.Size = 18
    $tr.Font.Underline = -1
    $tbl.Columns.Item(1).Width = EMU $widthEmu
    $tbl.Rows.Item(1).Height = EMU $heightEmu
}

# ---------------------------------------------------------------------
# 3) Table 16 (id=17) "ab0:AddressBook" -> "hms0:HotelManagementSystem"
#    Position unchanged, widened.
# ---------------------------------------------------------------------
$table16 = $s.Shapes.Item(3)
Set-StateTableLabel $table16 "hms0:HotelManagementSystem" 3207000 417888

# ---------------------------------------------------------------------
# 4) "Rectangle 19" - unchanged.
# 5) "TextBox 20" (currentStatePointer = 1, bottom row) - shift right.
# ---------------------------------------------------------------------
$textBox20 = $s.Shapes.Item(5)
$textBox20.Left = EMU 3909780
$textBox20.Top  = EMU 5170043

# ---------------------------------------------------------------------
# 6) "Rectangle 23" - unchanged.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 7) Table 14 (id=15) "ab2:AddressBook" -> "hms2:HotelManagementSystem"
#    Moved right, widened.
# ---------------------------------------------------------------------
$table14 = $s.Shapes.Item(7)
$table14.Left = EMU 7050408
$table14.Top  = EMU 1478952
Set-StateTableLabel $table14 "hms2:HotelManagementSystem" 3207000 417888

# ---------------------------------------------------------------------
# 8) Table 27 (id=28) "ab1:AddressBook" -> "hms1:HotelManagementSystem"
#    Moved right, widened.
# ---------------------------------------------------------------------
$table27 = $s.Shapes.Item(8)
$table27.Left = EMU 3761824
$table27.Top  = EMU 1476102
Set-StateTableLabel $table27 "hms1:HotelManagementSystem" 3207000 417888

# ---------------------------------------------------------------------
# 9) Table 28 (id=29) "ab0:AddressBook" -> "hms0:HotelManagementSystem"
#    Position unchanged, widened.
# ---------------------------------------------------------------------
$table28 = $s.Shapes.Item(9)
Set-StateTableLabel $table28 "hms0:HotelManagementSystem" 3207000 417888

# ---------------------------------------------------------------------
# 10) Table 29 (id=30) "ab2:AddressBook" -> "hms2:HotelManagementSystem"
#     Moved right, widened.
# ---------------------------------------------------------------------
$table29 = $s.Shapes.Item(10)
$table29.Left = EMU 7131994
$table29.Top  = EMU 3926589
Set-StateTableLabel $table29 "hms2:HotelManagementSystem" 3206997 417888

# ---------------------------------------------------------------------
# 11) Table 30 (id=31) "ab1:AddressBook" -> "hms1:HotelManagementSystem"
#     Moved right, widened.
# ---------------------------------------------------------------------
$table30 = $s.Shapes.Item(11)
$table30.Left = EMU 3761824
$table30.Top  = EMU 3926589
Set-StateTableLabel $table30 "hms1:HotelManagementSystem" 3288585 417888

# ---------------------------------------------------------------------
# 12) "Straight Arrow Connector 17" - shift right.
# ---------------------------------------------------------------------
$conn17 = $s.Shapes.Item(12)
$conn17.Left = EMU 5403681
$conn17.Top  = EMU 2024818

# ---------------------------------------------------------------------
# 13) "Straight Arrow Connector 21" - shift right.
# ---------------------------------------------------------------------
$conn21 = $s.Shapes.Item(13)
$conn21.Left = EMU 5367280
$conn21.Top  = EMU 4486056
